$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 386.1
$ws.Range("I33").Value = 386.1
$ws.Range("K33").Value = 386.1
$ws.Range("M33").Value = -157.1
$ws.Range("H69").Value = 1000
$ws.Range("I69").Value = 1000
$ws.Range("K69").Value = 3000
$ws.Range("M69").Value = -2126
$ws.Range("H72").Value = 1000
$ws.Range("I72").Value = 1000
$ws.Range("K72").Value = 9000
$ws.Range("M72").Value = -4632
$ws.Range("H92").Value = 2048
$ws.Range("I92").Value = 2048
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2048
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -800
$ws.Range("N92").ClearContents()
$ws.Range("H96").Value = 2026
$ws.Range("I96").Value = 2026
$ws.Range("K96").Value = 6078
$ws.Range("M96").Value = -4705
$ws.Range("H106").Value = 6744.5
$ws.Range("I106").Value = 6744.5
$ws.Range("K106").Value = 6744.5
$ws.Range("M106").Value = -6113.5
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 2926
$ws.Range("J137").Value = 5978
$ws.Range("L137").Value = 17934
$ws.Range("N137").Value = -23034
$ws.Range("H138").Value = 5104.8076
$ws.Range("I138").Value = 1311.9231
$ws.Range("J138").Value = 6369.1025
$ws.Range("K138").Value = 3935.7693
$ws.Range("L138").Value = 19107.3075
$ws.Range("M138").Value = 1204.2307
$ws.Range("N138").Value = -29387.3075
$ws.Range("H141").Value = 6798
$ws.Range("J141").Value = 2000
$ws.Range("L141").Value = 6000
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2031.8572
$ws.Range("I45").Value = 2099.25
$ws.Range("K45").Value = 2099.25
$ws.Range("M45").Value = -1722.25
$ws.Range("H74").Value = 2203.7
$ws.Range("I74").Value = 938.26666
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 938.26666
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -64.26666
$ws.Range("N74").Value = -7748
$ws.Range("H77").Value = 2203.7
$ws.Range("I77").Value = 938.26666
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 4691.3333
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -323.3333000000002
$ws.Range("N77").Value = -38736
$ws.Range("H122").Value = 2918.5217
$ws.Range("I122").Value = 2902.7368
$ws.Range("K122").Value = 8708.2104
$ws.Range("M122").Value = -6258.2104

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1041.6666
$ws.Range("I94").Value = 835.2
$ws.Range("J94").Value = 1299.75
$ws.Range("K94").Value = 835.2
$ws.Range("L94").Value = 1299.75
$ws.Range("M94").Value = -384.2
$ws.Range("N94").Value = -2201.75
$ws.Range("H105").Value = 2736.5
$ws.Range("I105").Value = 2736.5
$ws.Range("K105").Value = 2736.5
$ws.Range("M105").Value = -989.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6462.6665
$ws.Range("I31").Value = 2936.3333
$ws.Range("K31").Value = 2936.3333
$ws.Range("M31").Value = -2641.3333
$ws.Range("H34").Value = 6462.6665
$ws.Range("I34").Value = 2936.3333
$ws.Range("K34").Value = 2936.3333
$ws.Range("M34").Value = -2734.3333
$ws.Range("H132").Value = 3381.6428
$ws.Range("I132").Value = 2673.3333
$ws.Range("K132").Value = 8019.999899999999
$ws.Range("M132").Value = -5489.999899999999
$ws.Range("H134").Value = 3920.9167
$ws.Range("I134").Value = 3952.4443
$ws.Range("J134").Value = 3826.3333
$ws.Range("K134").Value = 11857.3329
$ws.Range("L134").Value = 11478.9999
$ws.Range("M134").Value = -9322.332900000001
$ws.Range("N134").Value = -16548.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 6172912
$ws.Range("J2").Value = 109
$ws.Range("L2").Value = 654
$ws.Range("N2").Value = -880
$ws.Range("H7").Value = 479
$ws.Range("I7").Value = 479
$ws.Range("K7").Value = 1437
$ws.Range("M7").Value = -1325
$ws.Range("H34").Value = 3115.4
$ws.Range("I34").Value = 1866.3334
$ws.Range("K34").Value = 5599.0002
$ws.Range("M34").Value = -5515.0002
$ws.Range("H93").Value = 1992
$ws.Range("I93").Value = 1992
$ws.Range("K93").Value = 5976
$ws.Range("M93").Value = -4104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8000
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730
$ws.Range("H73").Value = 8000
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064
$ws.Range("H102").Value = 3849.1875
$ws.Range("I102").Value = 3849.1875
$ws.Range("K102").Value = 3849.1875
$ws.Range("M102").Value = -2227.1875
$ws.Range("H113").Value = 2497
$ws.Range("I113").Value = 2489
$ws.Range("K113").Value = 2489
$ws.Range("M113").Value = -319

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1179.6
$ws.Range("I7").Value = 1179.6
$ws.Range("K7").Value = 1179.6
$ws.Range("M7").Value = -1067.6
$ws.Range("H40").Value = 4408.6665
$ws.Range("I40").Value = 4482.4287
$ws.Range("J40").Value = 4150.5
$ws.Range("K40").Value = 4482.4287
$ws.Range("L40").Value = 4150.5
$ws.Range("M40").Value = -4346.4287
$ws.Range("N40").Value = -4422.5
$ws.Range("H93").Value = 978.6667
$ws.Range("I93").Value = 978.6667
$ws.Range("K93").Value = 978.6667
$ws.Range("M93").Value = 269.3333
$ws.Range("H122").Value = 3001
$ws.Range("I122").Value = 3001
$ws.Range("K122").Value = 9003
$ws.Range("M122").Value = -6553
$ws.Range("H126").Value = 1179.6
$ws.Range("I126").Value = 1179.6
$ws.Range("K126").Value = 3538.8
$ws.Range("M126").Value = -1068.8
$ws.Range("H132").Value = 5000.6665
$ws.Range("I132").Value = 3004
$ws.Range("K132").Value = 9012
$ws.Range("M132").Value = -6482

Write-Output "Applied all changes"